# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") on Sheet1 holds the per-game strikeout totals. The saved
# stats were regenerated (K replaces the old Strike# figure), so every row's
# K value needs to be rewritten with the freshly computed number.
#
# Row 31 (game index 29) is unchanged by the regen (K was already 0), so it
# is intentionally left out of the table below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is (worksheet row, new K value for column G).
$kUpdates = @(
    @(2, 0),
    @(3, 0),
    @(4, 1),
    @(5, 0),
    @(6, 0),
    @(7, 0),
    @(8, 3),
    @(9, 2),
    @(10, 1),
    @(11, 1),
    @(12, 0),
    @(13, 1),
    @(14, 1),
    @(15, 0),
    @(16, 1),
    @(17, 1),
    @(18, 2),
    @(19, 0),
    @(20, 2),
    @(21, 1),
    @(22, 0),
    @(23, 0),
    @(24, 2),
    @(25, 0),
    @(26, 1),
    @(27, 1),
    @(28, 1),
    @(29, 1),
    @(30, 1),
    @(32, 1),
    @(33, 0),
    @(34, 1),
    @(35, 0),
    @(36, 0),
    @(37, 1),
    @(38, 1),
    @(39, 2),
    @(40, 1),
    @(41, 0),
    @(42, 1),
    @(43, 0),
    @(44, 1),
    @(45, 1),
    @(46, 1),
    @(47, 2),
    @(48, 2),
    @(49, 1),
    @(50, 3),
    @(51, 2),
    @(52, 2),
    @(53, 2),
    @(54, 1),
    @(55, 2)
)

foreach ($update in $kUpdates) {
    $row = $update[0]
    $newK = $update[1]
    # Column G is the 7th column ("K").
    $ws.Cells.Item($row, 7).Value = $newK
}
